$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lineup_Template")

# Fill in the lineup form with this week's submission
$ws.Range("B2").Value = "Mike K"
$ws.Range("B3").Value = 1
$ws.Range("B5").Value = "Sandy T"
$ws.Range("B6").Value = "Leah M"
$ws.Range("B7").Value = "Leah M"
$ws.Range("B8").Value = "Adam A"
$ws.Range("B9").Value = "John J"
$ws.Range("B10").Value = "Brooks K"
$ws.Range("B11").Value = "Bryson D"
$ws.Range("B12").Value = "Brooks K"
$ws.Range("B13").Value = "Brooks K"
$ws.Range("B14").Value = "Albert O"
$ws.Range("B15").Value = "Christy J"
$ws.Range("B16").Value = "Phil M"
$ws.Range("B17").Value = "Albert O"
$ws.Range("B18").Value = "Florida A/C"

# Normalize the redundant "filled" formatting on the CSR rows so it matches
# the plain bordered style used by the rest of the label column
$ws.Range("A14:A17").Interior.Pattern = -4142

$ws.Range("B18").Select()
